$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.021.47'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '2.947.14'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '375.55'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.40'
$ws.Range('E6').Value = '  -3.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.540'
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.47'
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0854'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '3.402.38'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.18'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.63'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '11.20'
$ws.Range('E16').Value = '  +49.79%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.940.60'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.00'
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('D19').Value = '51.008.97'
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('E20').Value = '  -6.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.51'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('E22').Value = '  -1.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '265.86'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.84'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('E25').Value = '  +6.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.15'
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.58'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.74'
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('E30').Value = '  -4.36%  '
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.90'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '33.54'
$ws.Range('E35').Value = '  -4.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0444'
$ws.Range('E36').Value = '  -2.24%  '
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.17'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.54'
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('E41').Value = '  -2.90%  '
$ws.Range('E42').Value = '  -4.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.60'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.41'
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  +2.49%  '
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.33'
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('D49').Value = '1.993.88'
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('E51').Value = '  +2.22%  '
